$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.506.57"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "3.546.54"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "585.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("E12").Value = "  -5.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").Value = "4.107.62"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "663.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +12.01%  "
$ws.Range("D16").Value = "69.637.83"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "3.537.92"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.969"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "3.786.51"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  -8.87%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "505.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.42%  "
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.137"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("E50").Value = "  +21.03%  "
$ws.Range("E51").Value = "  +61.50%  "
